$d = $word.ActiveDocument

# The "Requisitos" bullet list currently starts with the LOT2060 entry.
# The edit moves that entry (text + its trailing line break) from the
# start of the list to the end of the list (after the LOT2049 entry).

$lot2060Text = "LOT2060 -  Tecnologia de Biopolímeros  (Requisito)"
$lot2049Text = "LOT2049 -  Genética e Biotecnologia Vegetal  (Requisito)"

# Step 1: delete the LOT2060 entry (including its trailing manual line
# break) from its current location, by replacing it with nothing.
$deleted = $d.Content.Find.Execute("$lot2060Text^l", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $deleted) {
    throw "Could not find the LOT2060 requisito entry to remove."
}

# Step 2: locate the LOT2049 entry (including its trailing line break)
# and insert the LOT2060 entry, with its own line break, right after it.
$rng = $d.Content
$found = $rng.Find.Execute("$lot2049Text^l", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the LOT2049 requisito entry to insert after."
}

$rng.Collapse(0)
$rng.InsertAfter("$lot2060Text`v")
